$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Mutual Fund" -> shifts to D, etc.)
$ws.Range("C1").EntireColumn.Insert()

# New header for inserted column
$ws.Range("C1").Value = "Industry"

$industries = @(
    "Metals & Minerals Trading",
    "Power",
    "Finance",
    "Automobiles",
    "Power",
    "Insurance",
    "Banks",
    "Pharmaceuticals & Biotechnology",
    "Personal Products",
    "Banks",
    "Finance",
    "Banks",
    "IT - Software",
    "Power",
    "IT - Software",
    "Auto Components",
    "Food Products",
    "IT - Software",
    "Realty",
    "Retailing",
    "Power",
    "Banks",
    "Retailing",
    "Cement & Cement Products",
    "Ferrous Metals",
    "Realty"
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
